$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column L ("picHexCol") holds hex color strings like "#rrggbb".
# Swap the red and blue channels for every data row (2..361).
for ($r = 2; $r -le 361; $r++) {
    $cell = $ws.Cells.Item($r, 12)
    $val = $cell.Value2
    if ($val -and $val.Length -eq 7 -and $val.Substring(0,1) -eq "#") {
        $rr = $val.Substring(1,2)
        $gg = $val.Substring(3,2)
        $bb = $val.Substring(5,2)
        $cell.Value = "#" + $bb + $gg + $rr
    }
}

# Touch row 362 (column L, the last used column) so the sheet's used
# range / dimension extends to include a new, otherwise-empty row 362 -
# mirrors the appended blank row seen in the target workbook. Reset the
# formatting straight back to the default style so no visible change
# (and minimal style-table footprint) remains.
$tail = $ws.Cells.Item(362, 12)
$tail.WrapText = $false
$tail.Style = "Normal"
